$d = $word.ActiveDocument

# --- Op 1: "en petit, " -> "en petit, il" -------------------------------
$d.Content.Find.Execute("en petit, ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "en petit, il", 2)

# --- Op 2: turn the <corr>il</corr> run-trio into <corr><del>z</del</corr>
# Locate the (unique) "<corr>" run to anchor the following operations.
$probe = $d.Content
$probe.Find.Execute("<corr>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

# 2a. "<corr>" -> "<corr><del>"
$scoped = $d.Range($probe.Start, $d.Content.End)
$scoped.Find.Execute("<corr>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<corr><del>", 1)

# 2b. the following run ("il") -> "z", and drop its color (000000) so
#     that it matches the plain/no-color run formatting used elsewhere.
$afterCorr = $d.Range($scoped.End, $d.Content.End)
$afterCorr.Find.Execute("il", $true, $false, $false, $false, $false, `
    $true, 1, $false, "z", 1)
$afterCorr.Font.Color = -16777216   # wdColorAutomatic: clears explicit color

# 2c. the following run ("</corr>") -> "</del</corr>"
$afterZ = $d.Range($afterCorr.End, $d.Content.End)
$afterZ.Find.Execute("</corr>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "</del</corr>", 1)

# --- Op 3: delete the trailing " " run right after "</del>" (before the
# comment range start) -----------------------------------------------
$afterDelCorr = $d.Range($afterZ.End, $d.Content.End)
$afterDelCorr.Find.Execute("</del> ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "</del>", 1)

Write-Host "done"
